$wb = $excel.ActiveWorkbook

# --- Matches sheet updates ---
$wsMatches = $wb.Worksheets.Item("Matches")
$wsMatches.Activate()

$wsMatches.Range("G2").Value = 1
$wsMatches.Range("H2").Value = 1

$wsMatches.Range("G3").Value = 1
$wsMatches.Range("H3").Value = 1

$wsMatches.Range("G4").Value = 1
$wsMatches.Range("H4").Value = 1

$wsMatches.Range("E7").Select()

# --- Scorers sheet updates ---
$wsScorers = $wb.Worksheets.Item("Scorers")
$wsScorers.Activate()

$wsScorers.Range("C2").Value = 1
$wsScorers.Range("C3").Value = 1

$wsScorers.Range("D4").Select()

# Return to Matches tab as the active sheet (matches original workbook state)
$wsMatches.Activate()
